# Final-Project-Test-Cases-version-2.xlsx — add three new "length" test
# cases (rows 35-37) to the "main" worksheet, including a hyperlink that
# Excel auto-creates on the "Seneca@Ontario979" value, and update the
# saved view (zoom/scroll/selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# --- Row 35: Manav$tudent -------------------------------------------------
$ws.Range("C35").Value = 'Manav$tudent'
$ws.Range("D35").Value = '"The length is 12"'
$ws.Range("E35").Value = "SAME"
$ws.Range("F35").Value = "PASS"

# --- Row 36: Seneca@Ontario979 (becomes a mailto hyperlink) --------------
$ws.Range("C36").Value = 'Seneca@Ontario979'
$ws.Hyperlinks.Add($ws.Range("C36"), "mailto:Seneca@Ontario979")

# --- Row 37: CanadaToronto987 (note: D37 typed before C37/D36) -----------
$ws.Range("D37").Value = '"The length is 16"'
$ws.Range("E36").Value = "SAME"
$ws.Range("F36").Value = "PASS"
$ws.Range("C37").Value = 'CanadaToronto987'
$ws.Range("D36").Value = '"The length is 17"'
$ws.Range("E37").Value = "SAME"
$ws.Range("F37").Value = "PASS"

# --- Comments column for all three rows -----------------------------------
$ws.Range("G35").Value = "The return length is correct."
$ws.Range("G36").Value = "The return length is correct."
$ws.Range("G37").Value = "The return length is correct."

# --- Update the saved view: scroll position, zoom, selection -------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 78
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H37").Select()

$wb.Save()
